$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("puns")

# Fix row 6 (How do attorneys sleep) - swap answer/hint into correct columns
$ws.Range("B6").Value = "First they lie on one side, then they lie on the other side."
$ws.Range("C6").Value = "Think about the positions and how you'd sleep."

# Fix row 8 (How do you make holy water) - swap answer/hint into correct columns
$ws.Range("B8").Value = "You take some regular water and boil the hell out of it."
$ws.Range("C8").Value = "Same way if you really want to make overcooked rice."

# Add new row 76 with a new pun (question, hint, answer)
# Shared-string append order must be: question, hint, answer
$ws.Range("A76").Value = "What kind of fish you find in a hospital"
$ws.Range("C76").Value = "The last specialist you'd want to see in a hospital."
$ws.Range("B76").Value = "A sturgeon."

# Refresh the view (zoom level and selected cell) similar to the saved session
$ws.Activate()
$ws.Range("B48").Select()
$excel.ActiveWindow.Zoom = 80
